$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Agregando reglas salariales mensuales: rename code HRBCD -> HRBDC
$ws.Range("A4").Value = "HRBDC"

# Update the active selection on the sheet (as left by the editor)
$ws.Range("A5").Select()
